# Append 17 new data rows (TRAN_NUMBER / Email_ID pairs) to the single
# worksheet, right after the existing last row (317), extending the used
# range to A1:B334.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("3717170JU", "n3f5f2dp@hotmail.com"),
    @("7410748VL", "r18wvs0o@hotmail.com"),
    @("8321043YJ", "e919iny9@hotmail.com"),
    @("0453228CM", "dsuokk99@hotmail.com"),
    @("6516713SQ", "6jxtabx4@hotmail.com"),
    @("3434755IM", "mspq2uxt@hotmail.com"),
    @("4538473NP", "tvn6r5nu@hotmail.com"),
    @("2218587GH", "ikg6v75c@hotmail.com"),
    @("8401321YN", "v452elbb@hotmail.com"),
    @("5530257PQ", "wxmaju3k@hotmail.com"),
    @("4742405mv", "s5qjrawy@hotmail.com"),
    @("2186845gf", "o2pshm5u@hotmail.com"),
    @("7185388vd", "4f9vm964@hotmail.com"),
    @("0136522bf", "i38u5evj@hotmail.com"),
    @("2511750go", "jugf2xax@hotmail.com"),
    @("3531388ko", "b5hzh728@hotmail.com"),
    @("1355752dj", "fnvx4vl1@hotmail.com")
)

$startRow = 318
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
}
